# updating loading app, some performance changes
$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A1").NumberFormat = "@"
$ws1.Range("A1").Value = "0,627401"

$ws1.Range("B1").ClearContents()
$ws1.Range("C1").ClearContents()
$ws1.Range("D1").ClearContents()

$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "29466,8"

$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "83,67"

$ws1.Range("A5").NumberFormat = "@"
$ws1.Range("A5").Value = "0,221543"

# --- data sheet ---
$wsd = $wb.Worksheets.Item("data")

$wsd.Range("A1").Value = "ETHEREUM"
$wsd.Range("B1").Value = "BITCOIN"
$wsd.Range("E1").Value = "XRP"

$wsd.Range("E2").Value = "Sheet1"
$wsd.Range("F2").ClearContents()

$wsd.Range("B3").Value = "A2"
$wsd.Range("C3").Value = "A3"
$wsd.Range("E3").Value = "A1"
$wsd.Range("F3").ClearContents()

# --- Arkusz3 ---
$ws3 = $wb.Worksheets.Item("Arkusz3")
$ws3.Range("A1").NumberFormat = "@"
$ws3.Range("A1").Value = "0,627851"

# --- Tab / selection changes ---
$ws1.Activate()
$ws1.Range("K12").Select()
